$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 188, shifting existing rows 188..251 down to 189..252
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new record's data
$ws.Cells.Item(188, 1).Value = 5
$ws.Cells.Item(188, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(188, 3).Value = "Maule"
$ws.Cells.Item(188, 4).Value = 44559
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = 100114013
$ws.Cells.Item(188, 7).Value = "Zanahoria"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 400
$ws.Cells.Item(188, 11).Value = 7000
$ws.Cells.Item(188, 12).Value = 7000
$ws.Cells.Item(188, 13).Value = 7000
$ws.Cells.Item(188, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(188, 15).Value = "Región de Ñuble"
$ws.Cells.Item(188, 16).Value = 350
$ws.Cells.Item(188, 17).Value = 20
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same numeric date format used by the other date cells in column D
$ws.Cells.Item(188, 4).NumberFormat = $ws.Cells.Item(189, 4).NumberFormat
